# "Carga Demandada" workbook — update installed power for the "Motores
# Elétricos" line with the correct engine power, letting the demand-load
# formulas (D6, D7) recalculate from it, and leave the selection on the
# cell that was last edited (B7).
#
# Note: the workbook's mc:AlternateContent/x15ac:absPath (the folder the
# file was last saved from, e.g. ".../Sewage\" -> ".../Sewage III\") is
# metadata Excel stamps internally from the OS save path; it has no
# corresponding property on the Workbook/Application object model (in
# real Excel or here), so it cannot be set via COM automation and is left
# untouched by this script.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Potência Instalada (installed power, kW) for Motores Elétricos, row 6.
$ws.Range("B6").Value = 29.46

# Recalculate so Potência Demandada (D6, shared formula B*C) and the
# Carga Total Demandada sum (D7 = SUM(D2:D6)) pick up the new input.
$excel.Calculate()

# Reflect the last-edited cell in the sheet's saved selection.
$ws.Range("B7").Select()

$wb.Save()
